# Apply the style-table changes described by the reference diff:
#   - Heading4: drop the explicit run color override (C00000 -> inherited/auto)
#   - VerbatimChar (character style linked to SourceCode): shrink font size 11pt -> 9pt
#   - SourceCode (paragraph style, linked to VerbatimChar): shrink font size 11pt -> 9pt
# (Heading5 only had its rsid stamp touched in the source diff; rsid values are
#  opaque editing-session identifiers that aren't exposed anywhere on the Word
#  object model, so there is no user-facing/content change to make there.)

$d = $word.ActiveDocument

# 1. Heading 4 style: remove the hard-coded dark-red run color so it once again
#    inherits its color (wdColorAutomatic clears the explicit override).
$heading4 = $d.Styles("Heading4")
$heading4.Font.Color = -16777216

# 2. Verbatim Char (character style, linked to the SourceCode paragraph style):
#    reduce the font size from 11pt (sz 22) to 9pt (sz 18).
$verbatimChar = $d.Styles("VerbatimChar")
$verbatimChar.Font.Size = 9

# 3. Source Code paragraph style: reduce the font size from 11pt (sz 22) to
#    9pt (sz 18) to match the linked character style above.
$sourceCode = $d.Styles("SourceCode")
$sourceCode.Font.Size = 9
